# Generate Report for Handoff
#
# The localization row for "b92cd537-eadb-45fe-a32c-870100e4edfc.md" has
# finished its translation pass and is now ready to hand off: update its
# Status / Priority / handoff-datetime on each language sheet, and mirror
# the status + latest-handoff-datetime on the Overview roll-up sheet.
# Because the new status text ("Ready for handoff") is longer than the old
# one ("In Translation"), the affected columns also grow a bit wider.

$wb = $excel.ActiveWorkbook

# ---- zh-cn sheet (row 3 = b92cd537 entry) ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("E3").Value = "mt"
$wsZh.Range("H3").Value = "2016-08-28 08:15:15"
$wsZh.Columns.Item(3).ColumnWidth = 16.3

# ---- de-de sheet (row 3 = b92cd537 entry) ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("E3").Value = "mt"
$wsDe.Range("H3").Value = "2016-08-28 08:15:20"
$wsDe.Columns.Item(3).ColumnWidth = 16.3

# ---- Overview sheet (row 3 = b92cd537 entry) ----
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Range("E3").Value = "Ready for handoff"
$wsOv.Range("F3").Value = "Ready for handoff"
$wsOv.Range("G3").Value = "2016-08-28 08:15:20"
$wsOv.Columns.Item(5).ColumnWidth = 16.3
$wsOv.Columns.Item(6).ColumnWidth = 16.3
